# Fuel Prod Imp Exp Balancing Priorities.xlsx - "updated 4.0 files and mdl"
#
# Semantic edits captured by this script:
#  1. About!C1  - bump the "last updated" date from 1/3/2024 to 3/28/2024.
#  2. FPIEBP!B3:D3 (hard coal priorities) - re-ordered from
#     production=3, imports=2, exports=1  ->  production=1, imports=3, exports=2.
#  3. Selection left on FPIEBP!E3 (was F4) to match where the author's
#     cursor ended up when the workbook was saved.

$wb = $excel.ActiveWorkbook

# --- About sheet -----------------------------------------------------
$about = $wb.Worksheets.Item("About")
$about.Range("C1").Value = "3/28/2024"

# --- FPIEBP sheet ------------------------------------------------------
$fpiebp = $wb.Worksheets.Item("FPIEBP")
$fpiebp.Activate()

$fpiebp.Range("B3").Value = 1
$fpiebp.Range("C3").Value = 3
$fpiebp.Range("D3").Value = 2

# Leave the selection where the author left it before saving.
$fpiebp.Range("E3").Select() | Out-Null
